$d = $word.ActiveDocument

# Update the date header (unique text, safe to replace across the whole document)
$d.Content.Find.Execute("2025-04-02 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-04-03 Thursday", 2) | Out-Null

# Update the multiplication-table cells.
# Several cells share identical "NNNxN=" text, so each cell is addressed individually
# via Tables(1).Cell(row, col).Range, and Replace is restricted to a single match
# (wdReplaceOne = 1) scoped to that cell so it cannot bleed into other cells that
# still contain the same original text.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Find.Execute("483×2=", $true, $false, $false, $false, $false, $true, 1, $false, "880×6=", 1) | Out-Null
$t.Cell(1, 2).Range.Find.Execute("965×4=", $true, $false, $false, $false, $false, $true, 1, $false, "727×9=", 1) | Out-Null
$t.Cell(1, 3).Range.Find.Execute("686×9=", $true, $false, $false, $false, $false, $true, 1, $false, "771×5=", 1) | Out-Null
$t.Cell(1, 4).Range.Find.Execute("196×8=", $true, $false, $false, $false, $false, $true, 1, $false, "418×8=", 1) | Out-Null
$t.Cell(1, 5).Range.Find.Execute("774×9=", $true, $false, $false, $false, $false, $true, 1, $false, "415×5=", 1) | Out-Null
$t.Cell(5, 1).Range.Find.Execute("896×2=", $true, $false, $false, $false, $false, $true, 1, $false, "641×5=", 1) | Out-Null
$t.Cell(5, 2).Range.Find.Execute("790×6=", $true, $false, $false, $false, $false, $true, 1, $false, "581×2=", 1) | Out-Null
$t.Cell(5, 3).Range.Find.Execute("483×2=", $true, $false, $false, $false, $false, $true, 1, $false, "106×3=", 1) | Out-Null
$t.Cell(5, 4).Range.Find.Execute("870×9=", $true, $false, $false, $false, $false, $true, 1, $false, "521×4=", 1) | Out-Null
$t.Cell(5, 5).Range.Find.Execute("898×6=", $true, $false, $false, $false, $false, $true, 1, $false, "835×9=", 1) | Out-Null
$t.Cell(10, 1).Range.Find.Execute("461×4=", $true, $false, $false, $false, $false, $true, 1, $false, "690×8=", 1) | Out-Null
$t.Cell(10, 2).Range.Find.Execute("447×2=", $true, $false, $false, $false, $false, $true, 1, $false, "147×2=", 1) | Out-Null
$t.Cell(10, 3).Range.Find.Execute("954×6=", $true, $false, $false, $false, $false, $true, 1, $false, "535×4=", 1) | Out-Null
$t.Cell(10, 4).Range.Find.Execute("102×2=", $true, $false, $false, $false, $false, $true, 1, $false, "194×2=", 1) | Out-Null
$t.Cell(10, 5).Range.Find.Execute("447×2=", $true, $false, $false, $false, $false, $true, 1, $false, "405×5=", 1) | Out-Null
$t.Cell(15, 1).Range.Find.Execute("985×9=", $true, $false, $false, $false, $false, $true, 1, $false, "597×7=", 1) | Out-Null
$t.Cell(15, 2).Range.Find.Execute("293×9=", $true, $false, $false, $false, $false, $true, 1, $false, "920×7=", 1) | Out-Null
$t.Cell(15, 3).Range.Find.Execute("657×5=", $true, $false, $false, $false, $false, $true, 1, $false, "770×5=", 1) | Out-Null
$t.Cell(15, 4).Range.Find.Execute("531×4=", $true, $false, $false, $false, $false, $true, 1, $false, "656×9=", 1) | Out-Null
$t.Cell(15, 5).Range.Find.Execute("774×2=", $true, $false, $false, $false, $false, $true, 1, $false, "314×6=", 1) | Out-Null
$t.Cell(20, 1).Range.Find.Execute("823×6=", $true, $false, $false, $false, $false, $true, 1, $false, "804×7=", 1) | Out-Null
$t.Cell(20, 2).Range.Find.Execute("479×9=", $true, $false, $false, $false, $false, $true, 1, $false, "309×5=", 1) | Out-Null
$t.Cell(20, 3).Range.Find.Execute("652×2=", $true, $false, $false, $false, $false, $true, 1, $false, "355×8=", 1) | Out-Null
$t.Cell(20, 4).Range.Find.Execute("136×8=", $true, $false, $false, $false, $false, $true, 1, $false, "697×9=", 1) | Out-Null
$t.Cell(20, 5).Range.Find.Execute("522×3=", $true, $false, $false, $false, $false, $true, 1, $false, "216×9=", 1) | Out-Null

Write-Output "Done updating date header and 25 multiplication-table cells."
